$d = $word.ActiveDocument

# Locate the exact sentence that needs to change. We search for the full
# original tail of the paragraph so the match is unambiguous and so the
# found Range's Start/End give us precise character offsets to work with.
$oldText = ") or using MPEI package installer. Source code is available on Google Code site."

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target sentence to update."
}

$oldStart = $rng.Start
$oldEnd = $rng.End

# Build the replacement out of the same pieces used in the authored edit
# (kept separate here for clarity / traceability even though adjoining
# runs that end up with identical formatting are written back out as a
# single run).
$pieces = @(
    ")",
    ", using MPEI package installer or at our Google Code ",
    "site",
    ". ",
    "Source code is ",
    "also available on Google Code",
    "."
)
$newText = [string]::Join("", $pieces)

# Insert the new text right after the old text (so we don't disturb the
# formatting of the hyperlink run that immediately precedes it), then
# remove the old text. Doing it in this order avoids the replacement
# picking up the hyperlink's character style.
$insertionPoint = $d.Range($oldEnd, $oldEnd)
$insertionPoint.InsertBefore($newText)

$d.Range($oldStart, $oldEnd).Delete()
